# Commit: "integrate database all in one"
# Updates the latest-quarter ("12 ماهه منتهی به 1401/12", column M) figures
# and refreshes the two "published on" date labels (columns I and M, row 9)
# from 1402-03-13 to 1402-04-06 on the Overview sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 9: "تاریخ انتشار" (publish date) labels ---
# Column I (9 ماهه منتهی به 1401/09) publish-date tag
$ws.Cells.Item(9, 9).Value = "1402-04-06 (11)"
# Column M (12 ماهه منتهی به 1401/12) publish-date tag
$ws.Cells.Item(9, 13).Value = "1402-04-06 (3)"

# --- Column M (row 13 header = 12 ماهه منتهی به 1401/12) figures ---
$ws.Cells.Item(14, 13).Value = -9694342    # هزینه های عمومی, اداری و تشکیلاتی
$ws.Cells.Item(17, 13).Value = 12641483    # سود (زیان) عملیاتی
$ws.Cells.Item(19, 13).Value = 27335571    # سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Cells.Item(20, 13).Value = 32156798    # مالیات
$ws.Cells.Item(21, 13).Value = -706151     # سود (زیان) خالص عملیات در حال تداوم
$ws.Cells.Item(22, 13).Value = 31450647    # سود (زیان) عملیات متوقف شده پس از اثر مالیاتی
$ws.Cells.Item(24, 13).Value = 31450647    # سود (زیان) خالص
$ws.Cells.Item(25, 13).Value = 1016        # سود هر سهم پس از کسر مالیات
$ws.Cells.Item(27, 13).Value = 1016        # سود هر سهم بر اساس آخرین سرمایه
